# Fix AWS pricing in the benchmark XLS file.
# The source hourly instance prices (column B) dropped; every dependent
# formula cell (K:P on the rate rows, H:M on the per-job-cost rows) recalcs
# automatically once these inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated hourly instance prices (AWS price correction) ---
$ws.Range("B25").Value = 0.07    # was 0.113
$ws.Range("B26").Value = 0.14    # was 0.225
$ws.Range("B27").Value = 0.28    # was 0.45
$ws.Range("B28").Value = 0.56    # was 0.9

$ws.Range("B32").Value = 0.105   # was 0.15
$ws.Range("B33").Value = 0.21    # was 0.3
$ws.Range("B34").Value = 0.42    # was 0.6

# --- Row heights shrink slightly to match the new layout ---
$ws.Rows.Item(25).RowHeight = 12.65
$ws.Rows.Item(26).RowHeight = 12.65
$ws.Rows.Item(27).RowHeight = 12.65
$ws.Rows.Item(28).RowHeight = 12.65
$ws.Rows.Item(32).RowHeight = 12.65
$ws.Rows.Item(33).RowHeight = 12.65
$ws.Rows.Item(34).RowHeight = 12.65

# --- View state: scroll position, tab ratio, and current selection ---
$excel.ActiveWindow.TabRatio = 211
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F42").Select()
